# Insert a new data row at row 47 (pushing existing rows 47-54 down to 48-55)
# and populate it with the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(47).Insert()

$ws.Cells.Item(47, 1).Value = 5
$ws.Cells.Item(47, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(47, 3).Value = "Maule"
$ws.Cells.Item(47, 4).Value = 44468
$ws.Cells.Item(47, 5).Value = 7
$ws.Cells.Item(47, 6).Value = 100112013
$ws.Cells.Item(47, 7).Value = "Alcachofa"
$ws.Cells.Item(47, 8).Value = "Madrigal"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 300
$ws.Cells.Item(47, 11).Value = 10000
$ws.Cells.Item(47, 12).Value = 10000
$ws.Cells.Item(47, 13).Value = 10000
$ws.Cells.Item(47, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(47, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(47, 16).Value = 250
$ws.Cells.Item(47, 17).Value = 40
$ws.Cells.Item(47, 18).Value = "Hortaliza"
